$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.702.13'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.710.12'
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.55%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').Value = '2.706.94'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.361'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Value = '3.208.07'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '68.679.46'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').Value = '2.685.12'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('E22').Value = '  +2.60%  '
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '594.30'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.32%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('E33').Value = '  +2.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('E35').Value = '  +2.63%  '
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.91'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '160.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.381'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  -5.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '157.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('E48').Value = '  +5.28%  '
$ws.Range('E49').Value = '  +5.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.608'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.37%  '

Write-Output "Updated cryptos list"